# Add a new "Results" slide (Title and Content layout) at the end of the deck,
# mirroring the existing slide 7 ("Results") but with the body content filled in.

$p = $ppt.ActivePresentation

# ppLayoutText = 2 ("Title and Content") -> same layout used by the rest of the deck.
$s = $p.Slides.Add($p.Slides.Count + 1, 2)

# --- Title -----------------------------------------------------------------
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Results"

# --- Body content ------------------------------------------------------------
$tf = $s.Shapes.Item(2).TextFrame
$tr = $tf.TextRange

$tr.Text = "Significant speedup for large matrices`rStrassen"
[void]$tr.InsertAfter(" algorithm slower for smaller matrices")
[void]$tr.InsertAfter("`rParallel algorithm is faster than serial one.")
[void]$tr.InsertAfter("`rGranularity changes for different matrix sizes to avoid memory overflow.")
[void]$tr.InsertAfter("`r")
